$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = [double]"6.505707639270497E-08"
$ws.Range("C2").Value = [double]"5.351274978693255E-14"

$ws.Range("B3").Value = [double]"5.954015320330797E-06"
$ws.Range("C3").Value = [double]"1.998401444325282E-15"
$ws.Range("H3").Value = 0

$ws.Range("B4").Value = [double]"1.521301817852816E-08"
$ws.Range("C4").Value = [double]"8.060219158778636E-14"

$ws.Range("B5").Value = [double]"6.646301049251235E-05"
$ws.Range("C5").Value = [double]"2.220446049250313E-16"

$ws.Range("C7").Value = 0

$ws.Range("C8").Value = [double]"3.219224886663596E-11"
$ws.Range("E8").Value = 0
$ws.Range("H8").Value = 0

$ws.Range("B9").Value = [double]"2.395753306849713E-08"
$ws.Range("C9").Value = [double]"8.659739592076221E-15"

$ws.Range("B10").Value = 0.74812683126579
$ws.Range("C10").Value = 0.06208280264152122
$ws.Range("D10").Value = 0.1341533866968256
$ws.Range("E10").Value = 0.08882838250278091
$ws.Range("F10").Value = [double]"3.311586657472887E-05"
